# Weekly update: insert 3 new rows of data (new "Hass" quality grades from
# Cabildo, boxed at 17 kilos) at the top of the "Vega Monumental Concepcion -
# Palta" weekly block, pushing the existing rows 310-334 down to 313-337.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current row 310, shifting rows
# 310:334 down to 313:337 (keeps all of their existing data/formatting,
# including the date-format style on column D).
$ws.Range("A310:T312").EntireRow.Insert()

# New row 310 - Hass, "1a nueva(o)", Cabildo, caja de 17 kilos.
$ws.Cells.Item(310, 1).Value = 11
$ws.Cells.Item(310, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(310, 3).Value = "Bíobío"
$ws.Cells.Item(310, 4).Value = 44491
$ws.Cells.Item(310, 5).Value = 8
$ws.Cells.Item(310, 6).Value = "Fruta"
$ws.Cells.Item(310, 7).Value = 100106
$ws.Cells.Item(310, 8).Value = "Oleaginosos"
$ws.Cells.Item(310, 9).Value = 100106002
$ws.Cells.Item(310, 10).Value = "Palta"
$ws.Cells.Item(310, 11).Value = "Hass"
$ws.Cells.Item(310, 12).Value = "1a nueva(o)"
$ws.Cells.Item(310, 13).Value = 200
$ws.Cells.Item(310, 14).Value = 3000
$ws.Cells.Item(310, 15).Value = 3000
$ws.Cells.Item(310, 16).Value = 3000
$ws.Cells.Item(310, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(310, 18).Value = "Cabildo"
$ws.Cells.Item(310, 19).Value = 3000
$ws.Cells.Item(310, 20).Value = 1

# New row 311 - Hass, "2a nueva(o)", Cabildo, caja de 17 kilos.
$ws.Cells.Item(311, 1).Value = 11
$ws.Cells.Item(311, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(311, 3).Value = "Bíobío"
$ws.Cells.Item(311, 4).Value = 44491
$ws.Cells.Item(311, 5).Value = 8
$ws.Cells.Item(311, 6).Value = "Fruta"
$ws.Cells.Item(311, 7).Value = 100106
$ws.Cells.Item(311, 8).Value = "Oleaginosos"
$ws.Cells.Item(311, 9).Value = 100106002
$ws.Cells.Item(311, 10).Value = "Palta"
$ws.Cells.Item(311, 11).Value = "Hass"
$ws.Cells.Item(311, 12).Value = "2a nueva(o)"
$ws.Cells.Item(311, 13).Value = 250
$ws.Cells.Item(311, 14).Value = 2500
$ws.Cells.Item(311, 15).Value = 2500
$ws.Cells.Item(311, 16).Value = 2500
$ws.Cells.Item(311, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(311, 18).Value = "Cabildo"
$ws.Cells.Item(311, 19).Value = 2500
$ws.Cells.Item(311, 20).Value = 1

# New row 312 - Hass, "3a nueva (o)", Cabildo, caja de 17 kilos.
$ws.Cells.Item(312, 1).Value = 11
$ws.Cells.Item(312, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(312, 3).Value = "Bíobío"
$ws.Cells.Item(312, 4).Value = 44491
$ws.Cells.Item(312, 5).Value = 8
$ws.Cells.Item(312, 6).Value = "Fruta"
$ws.Cells.Item(312, 7).Value = 100106
$ws.Cells.Item(312, 8).Value = "Oleaginosos"
$ws.Cells.Item(312, 9).Value = 100106002
$ws.Cells.Item(312, 10).Value = "Palta"
$ws.Cells.Item(312, 11).Value = "Hass"
$ws.Cells.Item(312, 12).Value = "3a nueva (o)"
$ws.Cells.Item(312, 13).Value = 150
$ws.Cells.Item(312, 14).Value = 2000
$ws.Cells.Item(312, 15).Value = 2000
$ws.Cells.Item(312, 16).Value = 2000
$ws.Cells.Item(312, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(312, 18).Value = "Cabildo"
$ws.Cells.Item(312, 19).Value = 2000
$ws.Cells.Item(312, 20).Value = 1

# Column D carries a custom date/time number format in this sheet; make sure
# the freshly-inserted cells keep it explicitly (Insert already copies it
# down from the row above, but set it again to be safe).
$ws.Range("D310:D312").NumberFormat = "YYYY-MM-DD HH:MM:SS"
